$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 8773290
$ws.Range("I92").Value = 14493684
$ws.Range("K92").Value = 14493684
$ws.Range("M92").Value = -14492436
$ws.Range("H131").Value = 1915.2727
$ws.Range("I131").Value = 1625.625
$ws.Range("J131").Value = 2080.7856
$ws.Range("K131").Value = 4876.875
$ws.Range("L131").Value = 6242.3568
$ws.Range("M131").Value = 163.125
$ws.Range("N131").Value = -16322.3568
$ws.Range("H137").Value = 1210.2982
$ws.Range("I137").Value = 970.4103
$ws.Range("J137").Value = 1730.0555
$ws.Range("K137").Value = 2911.2309
$ws.Range("L137").Value = 5190.166499999999
$ws.Range("M137").Value = -361.2309
$ws.Range("N137").Value = -10290.1665
$ws.Range("H138").Value = 2740.54
$ws.Range("I138").Value = 1037.2354
$ws.Range("J138").Value = 3089.4097
$ws.Range("K138").Value = 3111.7062
$ws.Range("L138").Value = 9268.2291
$ws.Range("M138").Value = 2028.2938
$ws.Range("N138").Value = -19548.2291
$ws.Range("H141").Value = 5071.737
$ws.Range("I141").Value = 2264.9167
$ws.Range("J141").Value = 9883.429
$ws.Range("K141").Value = 6794.750100000001
$ws.Range("L141").Value = 29650.287
$ws.Range("M141").Value = -1614.750100000001
$ws.Range("N141").Value = -40010.287

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 484.3684
$ws.Range("I2").Value = 462.85715
$ws.Range("J2").Value = 544.6
$ws.Range("K2").Value = 462.85715
$ws.Range("L2").Value = 544.6
$ws.Range("M2").Value = -349.85715
$ws.Range("N2").Value = -770.6
$ws.Range("H63").Value = 29843.617
$ws.Range("I63").Value = 83541.63
$ws.Range("J63").Value = 4161.9565
$ws.Range("K63").Value = 83541.63
$ws.Range("L63").Value = 4161.9565
$ws.Range("M63").Value = -82855.63
$ws.Range("N63").Value = -5533.9565
$ws.Range("H66").Value = 29843.617
$ws.Range("I66").Value = 83541.63
$ws.Range("J66").Value = 4161.9565
$ws.Range("K66").Value = 417708.15
$ws.Range("L66").Value = 20809.7825
$ws.Range("M66").Value = -414276.15
$ws.Range("N66").Value = -27673.7825
$ws.Range("H110").Value = 180632.2
$ws.Range("I110").Value = 225577.75
$ws.Range("K110").Value = 225577.75
$ws.Range("M110").Value = -223532.75
$ws.Range("H116").Value = 484.3684
$ws.Range("I116").Value = 462.85715
$ws.Range("J116").Value = 544.6
$ws.Range("K116").Value = 462.85715
$ws.Range("L116").Value = 544.6
$ws.Range("M116").Value = 1831.14285
$ws.Range("N116").Value = -5132.6
$ws.Range("H132").Value = 5072.5264
$ws.Range("I132").Value = 4844
$ws.Range("J132").Value = 5567.6665
$ws.Range("K132").Value = 14532
$ws.Range("L132").Value = 16702.9995
$ws.Range("M132").Value = -12002
$ws.Range("N132").Value = -21762.9995

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 484.3684
$ws.Range("I3").Value = 462.85715
$ws.Range("J3").Value = 544.6
$ws.Range("K3").Value = 462.85715
$ws.Range("L3").Value = 544.6
$ws.Range("M3").Value = -348.85715
$ws.Range("N3").Value = -772.6
$ws.Range("H35").Value = 15333
$ws.Range("J35").Value = 15333
$ws.Range("L35").Value = 15333
$ws.Range("N35").Value = -15953
$ws.Range("H82").Value = 16935.834
$ws.Range("I82").Value = 5686.75
$ws.Range("J82").Value = 22560.375
$ws.Range("K82").Value = 5686.75
$ws.Range("L82").Value = 22560.375
$ws.Range("M82").Value = -5303.75
$ws.Range("N82").Value = -23326.375
$ws.Range("H85").Value = 16935.834
$ws.Range("I85").Value = 5686.75
$ws.Range("J85").Value = 22560.375
$ws.Range("K85").Value = 5686.75
$ws.Range("L85").Value = 22560.375
$ws.Range("M85").Value = -4360.75
$ws.Range("N85").Value = -25212.375
$ws.Range("H107").Value = 72668.71000000001
$ws.Range("I107").Value = 84280.164
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 84280.164
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -82360.164
$ws.Range("N107").Value = -6840
$ws.Range("H114").Value = 56325.332
$ws.Range("J114").Value = 56325.332
$ws.Range("L114").Value = 56325.332
$ws.Range("N114").Value = -65003.332

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1685.625
$ws.Range("I16").Value = 1721.5
$ws.Range("K16").Value = 1721.5
$ws.Range("M16").Value = -1434.5
$ws.Range("H31").Value = 5503.2856
$ws.Range("I31").Value = 1486.2903
$ws.Range("J31").Value = 8210.392
$ws.Range("K31").Value = 1486.2903
$ws.Range("L31").Value = 8210.392
$ws.Range("M31").Value = -1191.2903
$ws.Range("N31").Value = -8800.392
$ws.Range("H34").Value = 5503.2856
$ws.Range("I34").Value = 1486.2903
$ws.Range("J34").Value = 8210.392
$ws.Range("K34").Value = 1486.2903
$ws.Range("L34").Value = 8210.392
$ws.Range("M34").Value = -1284.2903
$ws.Range("N34").Value = -8614.392
$ws.Range("H58").Value = 1614.8889
$ws.Range("I58").Value = 1882.8
$ws.Range("J58").Value = 1511.8462
$ws.Range("K58").Value = 1882.8
$ws.Range("L58").Value = 1511.8462
$ws.Range("M58").Value = -1679.8
$ws.Range("N58").Value = -1917.8462
$ws.Range("H113").Value = 1685.625
$ws.Range("I113").Value = 1721.5
$ws.Range("K113").Value = 1721.5
$ws.Range("M113").Value = 448.5
$ws.Range("H136").Value = 1614.8889
$ws.Range("I136").Value = 1882.8
$ws.Range("J136").Value = 1511.8462
$ws.Range("K136").Value = 5648.4
$ws.Range("L136").Value = 4535.5386
$ws.Range("M136").Value = -3098.4
$ws.Range("N136").Value = -9635.5386

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 620.9863
$ws.Range("I5").Value = 484.19644
$ws.Range("J5").Value = 1071.5883
$ws.Range("K5").Value = 1452.58932
$ws.Range("L5").Value = 3214.7649
$ws.Range("M5").Value = -1340.58932
$ws.Range("N5").Value = -3438.7649
$ws.Range("H32").Value = 66670468
$ws.Range("J32").Value = 66670468
$ws.Range("L32").Value = 200011404
$ws.Range("N32").Value = -200011970
$ws.Range("H135").Value = 620.9863
$ws.Range("I135").Value = 484.19644
$ws.Range("J135").Value = 1071.5883
$ws.Range("K135").Value = 4357.76796
$ws.Range("L135").Value = 9644.294699999999
$ws.Range("M135").Value = -1822.76796
$ws.Range("N135").Value = -14714.2947
$ws.Range("H137").Value = 39357.4
$ws.Range("I137").Value = 7703.1113
$ws.Range("K137").Value = 23109.3339
$ws.Range("M137").Value = -18009.3339
$ws.Range("H141").Value = 9981.857
$ws.Range("I141").Value = 9202.111000000001
$ws.Range("K141").Value = 27606.333
$ws.Range("M141").Value = -22426.333

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 201.88889
$ws.Range("I2").Value = 201.88889
$ws.Range("K2").Value = 201.88889
$ws.Range("M2").Value = -88.88889
$ws.Range("H113").Value = 102039.37
$ws.Range("I113").Value = 124270.336
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 124270.336
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -122100.336
$ws.Range("N113").Value = -6340
$ws.Range("H132").Value = 2359.1086
$ws.Range("I132").Value = 1859.3928
$ws.Range("J132").Value = 3136.4443
$ws.Range("K132").Value = 5578.178400000001
$ws.Range("L132").Value = 9409.332900000001
$ws.Range("M132").Value = -3048.178400000001
$ws.Range("N132").Value = -14469.3329

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4140
$ws.Range("I96").Value = 3885.7144
$ws.Range("J96").Value = 4733.3335
$ws.Range("K96").Value = 3885.7144
$ws.Range("L96").Value = 4733.3335
$ws.Range("M96").Value = -2512.7144
$ws.Range("N96").Value = -7479.3335
$ws.Range("H107").Value = 715.9231
$ws.Range("I107").Value = 710.3333
$ws.Range("K107").Value = 2130.9999
$ws.Range("M107").Value = -210.9998999999998
$ws.Range("H122").Value = 2428.4
$ws.Range("I122").Value = 2428.4
$ws.Range("K122").Value = 7285.200000000001
$ws.Range("M122").Value = -4835.200000000001

Write-Host "Applied 202 cell updates across 7 worksheets."